$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.917.78"
$ws.Range("E2").Value = "  +5.94%  "
$ws.Range("D3").Value = "'3.146.94"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'590.59"
$ws.Range("E5").Value = "  +3.88%  "
$ws.Range("D6").Value = "'147.43"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'3.139.92"
$ws.Range("E8").Value = "  +3.90%  "
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("E10").Value = "  +17.88%  "
$ws.Range("D11").Value = "'5.73"
$ws.Range("E11").Value = "  +7.22%  "
$ws.Range("D12").Value = "'0.470"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "  +8.73%  "
$ws.Range("D14").Value = "'35.98"
$ws.Range("E14").Value = "  +4.92%  "
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "'3.667.65"
$ws.Range("E16").Value = "  +4.65%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'7.21"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'63.876.81"
$ws.Range("E18").Value = "  +5.94%  "
$ws.Range("D19").Value = "'3.142.31"
$ws.Range("E19").Value = "  +5.31%  "
$ws.Range("D20").Value = "'470.79"
$ws.Range("E20").Value = "  +5.89%  "
$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").Value = "'0.733"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E23").Value = "  +5.84%  "
$ws.Range("D24").Value = "'13.40"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'82.51"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'8.76"
$ws.Range("E27").Value = "  +10.38%  "
$ws.Range("D28").Value = "'2.71"
$ws.Range("E28").Value = "  +5.15%  "
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +7.65%  "
$ws.Range("D32").Value = "'27.09"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("D34").Value = "0.0₃0876"
$ws.Range("E34").Value = "  +9.20%  "
$ws.Range("E35").Value = "  +13.37%  "
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("E37").Value = "  +16.90%  "
$ws.Range("D38").Value = "'6.16"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").Value = "'50.92"
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("D40").Value = "'448.61"
$ws.Range("E40").Value = "  +10.44%  "
$ws.Range("D41").Value = "'8.74"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'0.0377"
$ws.Range("E42").Value = "  +5.88%  "
$ws.Range("D43").Value = "'2.925.82"
$ws.Range("E43").Value = "  +6.36%  "
$ws.Range("E44").Value = "  +10.30%  "
$ws.Range("E45").Value = "  +5.72%  "
$ws.Range("D46").Value = "'2.19"
$ws.Range("E46").Value = "  +6.70%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'34.92"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'123.98"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "'24.87"
$ws.Range("E51").Value = "  +5.33%  "
